$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the start-time ("F") values for rows 2-5
$ws.Range("F2").Value = 0.020833333333333332
$ws.Range("F3").Value = 0.010416666666666666
$ws.Range("F4").Value = 0.20833333333333334
$ws.Range("F5").Value = 0.020833333333333332

# Merge H2:H5, center the merged block, then give H2 the total elapsed
# time (formatted as a time value, same as column F/G)
$ws.Range("H2:H5").Merge()
$ws.Range("H2:H5").HorizontalAlignment = -4108  # xlCenter
$ws.Range("H2").NumberFormat = "h:mm"
$ws.Range("H2").Value = 0.2604166666666667

# Totals row
$ws.Range("F6").Formula = "=SUM(F2:F5)"
$ws.Range("G6").Formula = "=SUM(G2:G5)"
$ws.Range("H6").Formula = "=SUM(H2:H5)"

# Select H6, matching the final selection state in the workbook
$ws.Range("H6").Select()
